$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "노브랜드"
$ws.Range("B2").Value = "2024.04.30~05.08"
$ws.Range("C2").Value = "8,700~11,500"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = 10440
$ws.Range("F2").Value = "삼성증권"

$ws.Range("A3").Value = "아이씨티케이"
$ws.Range("B3").Value = "2024.04.24~04.30"
$ws.Range("C3").Value = "13,000~16,000"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 25610
$ws.Range("F3").Value = "NH투자증권"

$ws.Range("A4").Value = "이노그리드"
$ws.Range("B4").Value = "2024.04.18~04.24"
$ws.Range("C4").Value = "29,000~35,000"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = 17400
$ws.Range("F4").Value = "한국투자증권"

$ws.Range("A5").Value = "SK증권스팩12호"
$ws.Range("B5").Value = "2024.04.17~04.18"
$ws.Range("C5").Value = "2,000~2,000"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = 6000
$ws.Range("F5").Value = "SK증권"

$ws.Range("A6").Value = "HD현대마린솔루션(구.HD현대글로벌서비스)(유가)"
$ws.Range("B6").Value = "2024.04.16~04.22"
$ws.Range("C6").Value = "73,300~83,400"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = 652370
$ws.Range("F6").Value = "KB증권,신한투자증권,하나증권,대신증권,삼성증권"

$ws.Range("A7").Value = "코칩"
$ws.Range("B7").Value = "2024.04.15~04.19"
$ws.Range("C7").Value = "11,000~14,000"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = 16500
$ws.Range("F7").Value = "한국투자증권"

$ws.Range("A8").Value = "유안타스팩16호"
$ws.Range("B8").Value = "2024.04.15~04.16"
$ws.Range("C8").Value = "2,000~2,000"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = 10300
$ws.Range("F8").Value = "유안타증권"

$ws.Range("A9").Value = "민테크"
$ws.Range("B9").Value = "2024.04.12~04.18"
$ws.Range("C9").Value = "6,500~8,500"
$ws.Range("D9").Value = "-"
$ws.Range("E9").Value = 19500
$ws.Range("F9").Value = "KB증권"

$ws.Range("A10").Value = "디앤디파마텍"
$ws.Range("B10").Value = "2024.04.12~04.18"
$ws.Range("C10").Value = "22,000~26,000"
$ws.Range("D10").Value = "-"
$ws.Range("E10").Value = 24200
$ws.Range("F10").Value = "한국투자증권"

$ws.Range("A11").Value = "하나스팩33호"
$ws.Range("B11").Value = "2024.04.08~04.09"
$ws.Range("C11").Value = "2,000~2,000"
$ws.Range("D11").Value = "-"
$ws.Range("E11").Value = 7000
$ws.Range("F11").Value = "하나증권"

$ws.Range("A12").Value = "제일엠앤에스(구.제일기공)"
$ws.Range("B12").Value = "2024.04.05~04.12"
$ws.Range("C12").Value = "15,000~18,000"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = 36000
$ws.Range("F12").Value = "KB증권"

$ws.Range("A13").Value = "신한스팩13호"
$ws.Range("B13").Value = "2024.04.04~04.05"
$ws.Range("C13").Value = "2,000~2,000"
$ws.Range("D13").Value = "-"
$ws.Range("E13").Value = 6000
$ws.Range("F13").Value = "신한투자증권"

$ws.Range("A14").Value = "신한스팩12호"
$ws.Range("B14").Value = "2024.03.27~03.28"
$ws.Range("C14").Value = "2,000~2,000"
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = 10000
$ws.Range("F14").Value = "신한투자증권"

$ws.Range("A15").Value = "아이엠비디엑스"
$ws.Range("B15").Value = "2024.03.14~03.20"
$ws.Range("C15").Value = "7,700~9,900"
$ws.Range("D15").Value = "13000"
$ws.Range("E15").Value = 19250
$ws.Range("F15").Value = "미래에셋증권"

$ws.Range("A16").Value = "하나스팩32호"
$ws.Range("B16").Value = "2024.03.12~03.13"
$ws.Range("C16").Value = "2,000~2,000"
$ws.Range("D16").Value = "2000"
$ws.Range("E16").Value = 6000
$ws.Range("F16").Value = "하나증권"

$ws.Range("A17").Value = "엔젤로보틱스"
$ws.Range("B17").Value = "2024.03.06~03.12"
$ws.Range("C17").Value = "11,000~15,000"
$ws.Range("D17").Value = "20000"
$ws.Range("E17").Value = 17600
$ws.Range("F17").Value = "NH투자증권"

$ws.Range("A18").Value = "삼현"
$ws.Range("B18").Value = "2024.02.29~03.07"
$ws.Range("C18").Value = "20,000~25,000"
$ws.Range("D18").Value = "30000"
$ws.Range("E18").Value = 40000
$ws.Range("F18").Value = "한국투자증권"

$ws.Range("A19").Value = "오상헬스케어"
$ws.Range("B19").Value = "2024.02.21~02.27"
$ws.Range("C19").Value = "13,000~15,000"
$ws.Range("D19").Value = "20000"
$ws.Range("E19").Value = 12870
$ws.Range("F19").Value = "NH투자증권"

$ws.Range("A20").Value = "하나스팩31호"
$ws.Range("B20").Value = "2024.02.16~02.19"
$ws.Range("C20").Value = "2,000~2,000"
$ws.Range("D20").Value = "2000"
$ws.Range("E20").Value = 10000
$ws.Range("F20").Value = "하나증권"

$ws.Range("A21").Value = "케이엔알시스템"
$ws.Range("B21").Value = "2024.02.16~02.22"
$ws.Range("C21").Value = "9,000~11,000"
$ws.Range("D21").Value = "13500"
$ws.Range("E21").Value = 18936
$ws.Range("F21").Value = "DB금융투자,NH투자증권"
